# Insert a new "Player Info" worksheet. Worksheets.Add() with no arguments
# inserts the new sheet before the first existing sheet, which is exactly
# where we want it (before "ODI Batting").
$wb = $excel.ActiveWorkbook

$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

# Re-fetch the other sheets by name now that the sheet collection has
# changed, so references line up with the correct worksheets.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the bold / bordered / centered header style used by the other
# sheets ("ODI Batting", "ODI Bowling") in the workbook.
$header = $playerInfo.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# "5479" looks numeric, so force it to stay text (matching the original
# file's convention of storing every value, even numeric-looking ones, as
# text) instead of letting it get auto-converted into a number.
$idCell = $playerInfo.Range("A2")
$idCell.NumberFormat = "@"
$idCell.Value = "5479"
$idCell.ClearFormats()

$playerInfo.Range("B2").Value = "Bjorn Carl Fortuin"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Left Arm Orthodox"

# Update the "ODI Batting" sheet: rename MATCH_CARD_LINK -> MATCH_CODE
# and change the link values to bare match codes (kept as text, matching
# the original inline-string cell type, instead of being auto-converted
# to numbers).
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{ "D2" = "4405"; "D3" = "4657"; "D4" = "4658"; "D5" = "4727"; "D6" = "4731" }
foreach ($addr in $battingCodes.Keys) {
    $cell = $battingSheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $battingCodes[$addr]
    $cell.ClearFormats()
}

# Update the "ODI Bowling" sheet: rename MATCH_CARD_LINK -> MATCH_CODE
# and change the link values to bare match codes (kept as text).
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @{ "B2" = "4657"; "B3" = "4658"; "B4" = "4727"; "B5" = "4731" }
foreach ($addr in $bowlingCodes.Keys) {
    $cell = $bowlingSheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $bowlingCodes[$addr]
    $cell.ClearFormats()
}
